# major accuracy check update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the s2cDNAProtocol value (shared string) from "E7420" to "E7420L"
# for every data row (2-27). Assigning the whole range at once keeps all
# cells pointing at the same (updated) shared-string entry instead of
# minting a new one.
$ws.Range("G2:G27").Value = "E7420L"

# Replace the roboticS2Prep formulas ( =FALSE() ) with plain boolean
# literals (FALSE) for rows 2-27, preserving their existing style/format.
$ws.Range("H2:H27").Value = $false
